$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-04-09 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-04-10 Thursday", 2) | Out-Null

# Update the division problems in the table. Cell text is set directly
# (rather than via a global text search/replace) because several of the
# new values coincide with old values used elsewhere in the table, which
# would otherwise risk being replaced twice.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "16÷2="
$t.Cell(1,2).Range.Text  = "69÷8="
$t.Cell(1,3).Range.Text  = "98÷8="
$t.Cell(1,4).Range.Text  = "63÷3="
$t.Cell(1,5).Range.Text  = "79÷7="

$t.Cell(5,1).Range.Text  = "89÷8="
$t.Cell(5,2).Range.Text  = "75÷6="
$t.Cell(5,3).Range.Text  = "32÷9="
$t.Cell(5,4).Range.Text  = "32÷2="
$t.Cell(5,5).Range.Text  = "66÷5="

$t.Cell(9,1).Range.Text  = "52÷6="
$t.Cell(9,2).Range.Text  = "48÷5="
$t.Cell(9,3).Range.Text  = "59÷8="
$t.Cell(9,4).Range.Text  = "79÷6="
$t.Cell(9,5).Range.Text  = "45÷8="

$t.Cell(13,1).Range.Text = "84÷4="
$t.Cell(13,2).Range.Text = "17÷9="
$t.Cell(13,3).Range.Text = "15÷5="
$t.Cell(13,4).Range.Text = "58÷2="
$t.Cell(13,5).Range.Text = "70÷5="

$t.Cell(17,1).Range.Text = "42÷2="
$t.Cell(17,2).Range.Text = "95÷9="
$t.Cell(17,3).Range.Text = "43÷3="
$t.Cell(17,4).Range.Text = "20÷6="
$t.Cell(17,5).Range.Text = "73÷3="
